$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1721.653791726028
$ws.Range("D2").Value = 0.3668478529267383
$ws.Range("E2").Value = 1159.563705625653

$ws.Range("C3").Value = 2493.062071715172
$ws.Range("D3").Value = 0.6571974450906863
$ws.Range("E3").Value = 2077.325241671097

$ws.Range("C4").Value = 3268.023409288708
$ws.Range("D4").Value = 0.9549276024599289
$ws.Range("E4").Value = 3018.415892156498

$ws.Range("C5").Value = 4052.064782911247
$ws.Range("D5").Value = 1.299146985731944
$ws.Range("E5").Value = 4106.453617917246

$ws.Range("C6").Value = 4846.765329286962
$ws.Range("D6").Value = 1.638351812457543
$ws.Range("E6").Value = 5178.640909434644
